# Segunda version de la pregunta problema
# Rewrites the first two paragraphs (title + body) into a title paragraph
# (with a trailing run) plus two body paragraphs, matching the target
# canonical OOXML produced by the authoring tool.

$d = $word.ActiveDocument

$targetXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordml" wp14:textId="643FFB9B">
  <w:pPr>
    <w:spacing w:after="160" w:line="259" w:lineRule="auto"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/>
      <w:b w:val="1"/>
      <w:bCs w:val="1"/>
      <w:i w:val="0"/>
      <w:iCs w:val="0"/>
      <w:caps w:val="0"/>
      <w:smallCaps w:val="0"/>
      <w:noProof w:val="0"/>
      <w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">Pregunta problema </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/>
      <w:noProof w:val="0"/>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordml" wp14:textId="03B1DE57">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:i w:val="0"/>
      <w:iCs w:val="0"/>
      <w:caps w:val="0"/>
      <w:smallCaps w:val="0"/>
      <w:noProof w:val="0"/>
      <w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:i w:val="0"/>
      <w:iCs w:val="0"/>
      <w:caps w:val="0"/>
      <w:smallCaps w:val="0"/>
      <w:noProof w:val="0"/>
      <w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>El restaurante ha experimentado una serie de problemas en la gestión de pedidos que han tenido consecuencias negativas para su funcionamiento. La falta de organización y control en el proceso de toma de pedidos ha afectado tanto al personal encargado de esta tarea como a los clientes.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordml" wp14:textId="0FC6D981">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:i w:val="0"/>
      <w:iCs w:val="0"/>
      <w:caps w:val="0"/>
      <w:smallCaps w:val="0"/>
      <w:noProof w:val="0"/>
      <w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:i w:val="0"/>
      <w:iCs w:val="0"/>
      <w:caps w:val="0"/>
      <w:smallCaps w:val="0"/>
      <w:noProof w:val="0"/>
      <w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>En el caso del personal, la mala gestión de los pedidos ha generado una sobrecarga de trabajo y estrés emocional en el equipo encargado de esta tarea. La falta de recursos y herramientas para llevar a cabo su trabajo de manera efectiva ha sido una de las principales causas de estos problemas.</w:t>
  </w:r>
</w:p>
'@

# Replace the full original content (title paragraph + single body paragraph)
# with the three new paragraphs described by $targetXml. Selecting the whole
# document range (including the trailing paragraph mark) and calling
# InsertXML with multiple <w:p> elements swaps both the run contents AND the
# paragraph properties (pPr) for every supplied paragraph.
$full = $d.Range(0, $d.Content.End)
[void]$full.InsertXML($targetXml)

# The InsertXML paragraph-properties import only honors the <w:rPr> child of
# <w:pPr> (paragraph-mark run formatting); other paragraph formatting such as
# <w:spacing> must be (re)applied through the native paragraph API.
$p1 = $d.Paragraphs(1)
$p1.SpaceAfter = 8
$p1.LineSpacingRule = 5
$p1.LineSpacing = 12.95
